# Update latest output (run 132)
$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("E3").Value = -39.98787974999996
$wsSchedule.Range("F3").Value = -0.881567013888888

# --- Detailed sheet updates ---
$wsDetailed.Range("B37").Value = 4.39442
$wsDetailed.Range("B38").Value = 35.14435
$wsDetailed.Range("B39").Value = 45.7518
$wsDetailed.Range("C39").Value = "historical"
$wsDetailed.Range("B40").Value = 57.03529
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 58.40626
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 58.48021
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 57.06
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 56.98
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 57.06
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("B46").Value = 56.98
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 57.4513
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value = 58.28671
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 59.20372
$wsDetailed.Range("B50").Value = 57.31
$wsDetailed.Range("B51").Value = 62.62125
$wsDetailed.Range("B52").Value = 60.43291
$wsDetailed.Range("B53").Value = 57.31
$wsDetailed.Range("B54").Value = 57.31
$wsDetailed.Range("B55").Value = 63.04179
$wsDetailed.Range("B56").Value = 65
$wsDetailed.Range("B57").Value = 64.7062
$wsDetailed.Range("B58").Value = 65
$wsDetailed.Range("B59").Value = 64.42498999999999
$wsDetailed.Range("B60").Value = 64.78793
$wsDetailed.Range("B61").Value = 65
$wsDetailed.Range("B62").Value = 63.25165
$wsDetailed.Range("B65").Value = 0.12132
$wsDetailed.Range("B66").Value = -5.85701
$wsDetailed.Range("B67").Value = -6.55354
$wsDetailed.Range("B68").Value = -7.28974
$wsDetailed.Range("B69").Value = -15.3603
$wsDetailed.Range("B70").Value = -16.71922
$wsDetailed.Range("B71").Value = -16.06576
$wsDetailed.Range("B72").Value = -23.5
$wsDetailed.Range("B73").Value = -23.5
$wsDetailed.Range("B74").Value = -15.83057
$wsDetailed.Range("B75").Value = -19.27294
$wsDetailed.Range("B76").Value = -23.5
$wsDetailed.Range("B77").Value = -23.5
$wsDetailed.Range("B78").Value = -13.5
$wsDetailed.Range("B79").Value = -7.33502
$wsDetailed.Range("B80").Value = -5.96132
$wsDetailed.Range("B81").Value = -6
$wsDetailed.Range("B82").Value = -4.96497
$wsDetailed.Range("B83").Value = -5.01
$wsDetailed.Range("B84").Value = 0.51
$wsDetailed.Range("B85").Value = 48.66757
$wsDetailed.Range("B86").Value = 56.28829
$wsDetailed.Range("B87").Value = 84.79000000000001
$wsDetailed.Range("B88").Value = 109.09804
$wsDetailed.Range("B89").Value = 180.82816
$wsDetailed.Range("B90").Value = 169.07583
$wsDetailed.Range("B91").Value = 143.60463
$wsDetailed.Range("B92").Value = 138.63513
$wsDetailed.Range("B93").Value = 128.9903
$wsDetailed.Range("B94").Value = 108.89
$wsDetailed.Range("B95").Value = 123.32092
$wsDetailed.Range("B96").Value = 108.89
$wsDetailed.Range("B97").Value = 105.79
